$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New test case written for IntegrationPointRDO (row 5):
# both Total Test Cases and Automated Test Cases increase by 1.
$ws.Range("B5").Value = 18
$ws.Range("C5").Value = 18

# Update the active selection to D5 as recorded in the sheet view.
$ws.Activate()
$ws.Range("D5").Select()
